$d = $word.ActiveDocument

$pairs = @(
    @("21×84=1764", "53×63=3339"),
    @("47×30=1410", "62×46=2852"),
    @("45×49=2205", "89×96=8544"),
    @("47×51=2397", "67×82=5494"),
    @("97×56=5432", "44×29=1276"),
    @("55×36=1980", "27×21=567"),
    @("40×65=2600", "35×55=1925"),
    @("16×33=528", "49×60=2940"),
    @("54×83=4482", "97×76=7372"),
    @("28×79=2212", "94×19=1786"),
    @("14×60=840", "53×33=1749"),
    @("51×69=3519", "47×71=3337"),
    @("37×25=925", "75×87=6525"),
    @("87×24=2088", "52×40=2080"),
    @("28×38=1064", "47×94=4418"),
    @("28×45=1260", "34×86=2924"),
    @("29×32=928", "57×26=1482"),
    @("58×26=1508", "89×80=7120"),
    @("33×34=1122", "81×79=6399"),
    @("27×17=459", "53×63=3339"),
    @("88×95=8360", "81×11=891"),
    @("86×51=4386", "18×33=594"),
    @("64×71=4544", "36×92=3312"),
    @("77×20=1540", "81×98=7938"),
    @("83×94=7802", "36×80=2880")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
